# MDT WIP procID -> ProcIdent
# On sheet "EenT gegevens", cell D2 holds the column-identifier label
# "ProcID" (row 1 holds the human label "recht"; row 2 holds the
# internal identifier). E2 (=D2) and the shared-formula range F2:R2
# (=E2) pick this value up automatically on recalculation.
# Also nudge the saved selection from D2 to D3, matching the workbook
# state captured after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EenT gegevens")

$ws.Range("D2").Value = "ProcIdent"

$ws.Activate()
$ws.Range("D3").Select()
